$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-valued columns (B, C, E) are written verbatim.
# Column D needs explicit Text number-format while assigning, since several
# values look numeric (e.g. "1.000", "98.00", "0.000006567") and would
# otherwise be silently re-parsed/rounded by Excel's numeric auto-detection.

# Row 2
$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '24.911.03'
$c.Style = "Normal"
$ws.Range("E2").Value = '  -4.13%  '

# Row 3
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '1.635.45'
$c.Style = "Normal"
$ws.Range("E3").Value = '  -6.34%  '

# Row 4
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '0.9985'
$c.Style = "Normal"
$ws.Range("E4").Value = '  -0.15%  '

# Row 5
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '232.36'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -6.58%  '

# Row 6
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '1.000'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +0.07%  '

# Row 7
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.4769'
$c.Style = "Normal"
$ws.Range("E7").Value = '  -6.21%  '

# Row 8
$ws.Range("E8").Value = '  -3.79%  '

# Row 9
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.2582'
$c.Style = "Normal"
$ws.Range("E9").Value = '  -6.62%  '

# Row 10
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.06101'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -1.48%  '

# Row 11
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.07009'
$c.Style = "Normal"
$ws.Range("E11").Value = '  -3.58%  '

# Row 12
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '1.635.84'
$c.Style = "Normal"
$ws.Range("E12").Value = '  -6.36%  '

# Row 13
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '14.52'
$c.Style = "Normal"
$ws.Range("E13").Value = '  -4.55%  '

# Row 14
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '0.5976'
$c.Style = "Normal"
$ws.Range("E14").Value = '  -8.65%  '

# Row 15
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '4.355'
$c.Style = "Normal"
$ws.Range("E15").Value = '  -7.03%  '

# Row 16
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '73.46'
$c.Style = "Normal"
$ws.Range("E16").Value = '  -5.51%  '

# Row 17
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '0.9999'
$c.Style = "Normal"
$ws.Range("E17").Value = '  -0.05%  '

# Row 18
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '0.9994'
$c.Style = "Normal"
$ws.Range("E18").Value = '  -0.05%  '

# Row 19
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '24.918.06'
$c.Style = "Normal"
$ws.Range("E19").Value = '  -4.12%  '

# Row 20
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '0.000006567'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -4.62%  '

# Row 21
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '11.18'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -5.95%  '

# Row 22
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '1.848.11'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -6.21%  '

# Row 23
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '4.355'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -2.81%  '

# Row 24
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '8.574'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -2.14%  '

# Row 25
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '5.259'
$c.Style = "Normal"
$ws.Range("E25").Value = '  -2.36%  '

# Row 26
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '133.26'
$c.Style = "Normal"
$ws.Range("E26").Value = '  -2.27%  '

# Row 27
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '14.80'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -3.32%  '

# Row 28
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '1.388'
$c.Style = "Normal"
$ws.Range("E28").Value = '  -8.67%  '

# Row 29
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '103.33'
$c.Style = "Normal"
$ws.Range("E29").Value = '  -2.35%  '

# Row 30
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '1.637'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -8.42%  '

# Row 31
$ws.Range("E31").Value = '  +1.86%  '

# Row 32
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '0.07683'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -6.55%  '

# Row 33
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '3.535'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -3.27%  '

# Row 34
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '0.9989'
$c.Style = "Normal"
$ws.Range("E34").Value = '  +0.01%  '

# Row 35
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '0.04297'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -7.88%  '

# Row 36
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '2.582'
$c.Style = "Normal"
$ws.Range("E36").Value = '  -2.65%  '

# Row 37
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '0.9233'
$c.Style = "Normal"
$ws.Range("E37").Value = '  -7.54%  '

# Row 38
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.5803'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -5.57%  '

# Row 39
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '2.543'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -8.57%  '

# Row 40
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.01525'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -5.70%  '

# Row 41
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.9994'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -0.01%  '

# Row 42
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.8097'
$c.Style = "Normal"
$ws.Range("E42").Value = '  +5.42%  '

# Row 43
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '98.00'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -2.57%  '

# Row 44
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '1.755'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -9.33%  '

# Row 45
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.3681'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -6.38%  '

# Row 46
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '4.687'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -6.60%  '

# Row 47
$ws.Range("B47").Value = 'Cronos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '0.05192'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -2.85%  '

# Row 48
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '0.1086'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -6.08%  '

# Row 49
$ws.Range("B49").Value = 'Aptos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '6.040'
$c.Style = "Normal"
$ws.Range("E49").Value = '  -5.02%  '

# Row 50
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '29.34'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -4.41%  '

# Row 51
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '0.9993'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -0.27%  '
